$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = '0 -> 5 -> 4 -> 1 -> 3 -> 2 -> 0'
$ws.Cells.Item(2, 3).Value = '[(0, -12), (4, -4), (6, 5), (-4, -4), (4, 12), (-4, 12), (0, -12)]'
$ws.Cells.Item(2, 4).Value = 81.83703435555695

$ws.Cells.Item(3, 2).Value = '0 -> 2 -> 1 -> 4 -> 3 -> 2 -> 1 -> 5 -> 0'
$ws.Cells.Item(3, 3).Value = '[(0, -12), (-4, 12), (-4, -4), (6, 5), (4, 12), (-4, 12), (-4, -4), (4, -4), (0, -12)]'
$ws.Cells.Item(3, 4).Value = 102.0090559675463

$ws.Cells.Item(4, 2).Value = '0 -> 1 -> 3 -> 5 -> 4 -> 3 -> 2 -> 0'
$ws.Cells.Item(4, 3).Value = '[(0, -12), (-4, -4), (4, 12), (4, -4), (6, 5), (4, 12), (-4, 12), (0, -12)]'
$ws.Cells.Item(4, 4).Value = 91.66352019776375

$ws.Cells.Item(5, 2).Value = '0 -> 5 -> 4 -> 3 -> 1 -> 2 -> 0'
$ws.Cells.Item(5, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, -4), (-4, 12), (0, -12)]'
$ws.Cells.Item(5, 4).Value = 83.66352019776376

$ws.Cells.Item(6, 2).Value = '0 -> 5 -> 4 -> 3 -> 1 -> 2 -> 0'
$ws.Cells.Item(6, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, -4), (-4, 12), (0, -12)]'
$ws.Cells.Item(6, 4).Value = 83.66352019776376

$ws.Cells.Item(7, 2).Value = '0 -> 5 -> 4 -> 3 -> 2 -> 1 -> 0'
$ws.Cells.Item(7, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(7, 4).Value = 58.38819816657173

$ws.Cells.Item(8, 2).Value = '0 -> 5 -> 4 -> 3 -> 2 -> 1 -> 0'
$ws.Cells.Item(8, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(8, 4).Value = 58.38819816657173

$ws.Cells.Item(9, 2).Value = '0 -> 5 -> 4 -> 3 -> 2 -> 1 -> 0'
$ws.Cells.Item(9, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(9, 4).Value = 58.38819816657173

$ws.Cells.Item(10, 2).Value = '0 -> 5 -> 4 -> 3 -> 2 -> 1 -> 0'
$ws.Cells.Item(10, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(10, 4).Value = 58.38819816657173

$ws.Cells.Item(11, 2).Value = '0 -> 5 -> 4 -> 3 -> 2 -> 1 -> 0'
$ws.Cells.Item(11, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(11, 4).Value = 58.38819816657173

$ws.Cells.Item(12, 2).Value = '0 -> 5 -> 3 -> 2 -> 4 -> 3 -> 2 -> 1 -> 0'
$ws.Cells.Item(12, 3).Value = '[(0, -12), (4, -4), (4, 12), (-4, 12), (6, 5), (4, 12), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(12, 4).Value = 85.37520932501253

$ws.Cells.Item(13, 2).Value = '0 -> 4 -> 3 -> 2 -> 1 -> 3 -> 2 -> 5 -> 0'
$ws.Cells.Item(13, 3).Value = '[(0, -12), (6, 5), (4, 12), (-4, 12), (-4, -4), (4, 12), (-4, 12), (4, -4), (0, -12)]'
$ws.Cells.Item(13, 4).Value = 102.0292258165963

$ws.Cells.Item(14, 2).Value = '0 -> 5 -> 4 -> 3 -> 2 -> 1 -> 0'
$ws.Cells.Item(14, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(14, 4).Value = 58.38819816657173

$ws.Cells.Item(15, 2).Value = '0 -> 5 -> 4 -> 3 -> 2 -> 1 -> 0'
$ws.Cells.Item(15, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(15, 4).Value = 58.38819816657173

$ws.Cells.Item(16, 2).Value = '0 -> 5 -> 4 -> 3 -> 2 -> 1 -> 0'
$ws.Cells.Item(16, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(16, 4).Value = 58.38819816657173

$ws.Cells.Item(17, 2).Value = '0 -> 5 -> 4 -> 1 -> 3 -> 2 -> 0'
$ws.Cells.Item(17, 3).Value = '[(0, -12), (4, -4), (6, 5), (-4, -4), (4, 12), (-4, 12), (0, -12)]'
$ws.Cells.Item(17, 4).Value = 81.83703435555695

$ws.Cells.Item(18, 2).Value = '0 -> 5 -> 4 -> 3 -> 2 -> 1 -> 0'
$ws.Cells.Item(18, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(18, 4).Value = 58.38819816657173

$ws.Cells.Item(19, 2).Value = '0 -> 5 -> 4 -> 3 -> 2 -> 1 -> 0'
$ws.Cells.Item(19, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(19, 4).Value = 58.38819816657173

$ws.Cells.Item(20, 2).Value = '0 -> 5 -> 3 -> 2 -> 1 -> 5 -> 4 -> 0'
$ws.Cells.Item(20, 3).Value = '[(0, -12), (4, -4), (4, 12), (-4, 12), (-4, -4), (4, -4), (6, 5), (0, -12)]'
$ws.Cells.Item(20, 4).Value = 84.19157274461199

$ws.Cells.Item(21, 2).Value = '0 -> 5 -> 4 -> 3 -> 2 -> 1 -> 0'
$ws.Cells.Item(21, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(21, 4).Value = 58.38819816657173

$ws.Cells.Item(22, 2).Value = '0 -> 5 -> 4 -> 3 -> 5 -> 2 -> 1 -> 0'
$ws.Cells.Item(22, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (4, -4), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(22, 4).Value = 84.27674198657004

$ws.Cells.Item(23, 2).Value = '0 -> 5 -> 4 -> 3 -> 2 -> 1 -> 0'
$ws.Cells.Item(23, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(23, 4).Value = 58.38819816657173

$ws.Cells.Item(24, 2).Value = '0 -> 5 -> 4 -> 3 -> 2 -> 1 -> 0'
$ws.Cells.Item(24, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(24, 4).Value = 58.38819816657173

$ws.Cells.Item(25, 2).Value = '0 -> 5 -> 4 -> 3 -> 2 -> 1 -> 0'
$ws.Cells.Item(25, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(25, 4).Value = 58.38819816657173

$ws.Cells.Item(26, 2).Value = '0 -> 5 -> 4 -> 3 -> 2 -> 4 -> 3 -> 1 -> 0'
$ws.Cells.Item(26, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, 12), (6, 5), (4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(26, 4).Value = 79.76340749158426

$ws.Cells.Item(27, 2).Value = '0 -> 5 -> 4 -> 1 -> 3 -> 2 -> 0'
$ws.Cells.Item(27, 3).Value = '[(0, -12), (4, -4), (6, 5), (-4, -4), (4, 12), (-4, 12), (0, -12)]'
$ws.Cells.Item(27, 4).Value = 81.83703435555695

$ws.Cells.Item(28, 2).Value = '0 -> 5 -> 4 -> 3 -> 2 -> 1 -> 0'
$ws.Cells.Item(28, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(28, 4).Value = 58.38819816657173

$ws.Cells.Item(29, 2).Value = '0 -> 5 -> 4 -> 3 -> 2 -> 1 -> 0'
$ws.Cells.Item(29, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(29, 4).Value = 58.38819816657173

$ws.Cells.Item(30, 2).Value = '0 -> 5 -> 4 -> 3 -> 2 -> 1 -> 0'
$ws.Cells.Item(30, 3).Value = '[(0, -12), (4, -4), (6, 5), (4, 12), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(30, 4).Value = 58.38819816657173

$ws.Cells.Item(31, 2).Value = '0 -> 5 -> 1 -> 4 -> 3 -> 2 -> 0'
$ws.Cells.Item(31, 3).Value = '[(0, -12), (4, -4), (-4, -4), (6, 5), (4, 12), (-4, 12), (-4, -4), (0, -12)]'
$ws.Cells.Item(31, 4).Value = 70.00905596754626
